# Update "想去人数" (number of people interested) figures that changed
# between the two site-generation runs.
#
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types, which
# aggregates the exhibition rows among others) both contain the same
# six events, so the same six F-column values need bumping in both
# sheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# Row -> (old, new) value map for each sheet, keyed by the F-column
# cell address holding the "想去人数" figure for that event.
$exhibitionUpdates = @{
    "F3"  = 5
    "F5"  = 77
    "F6"  = 5227
    "F8"  = 80
    "F10" = 350
    "F11" = 10
}

$allTypesUpdates = @{
    "F4"  = 5
    "F9"  = 77
    "F10" = 5227
    "F12" = 80
    "F15" = 350
    "F16" = 10
}

foreach ($cellAddr in $exhibitionUpdates.Keys) {
    $wsExhibition.Range($cellAddr).Value = $exhibitionUpdates[$cellAddr]
}

foreach ($cellAddr in $allTypesUpdates.Keys) {
    $wsAll.Range($cellAddr).Value = $allTypesUpdates[$cellAddr]
}

$wb.Save()
